$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, shifting existing rows 8-54 down to 9-55
$ws.Rows.Item(8).Insert()

# Populate the new row 8 with the new data record
$ws.Range("A8").Value = 11
$ws.Range("B8").Value = "Vega Monumental Concepción"
$ws.Range("C8").Value = "Bíobío"
$ws.Range("D8").Value = 44764
$ws.Range("E8").Value = 8
$ws.Range("F8").Value = 100112013
$ws.Range("G8").Value = "Alcachofa"
$ws.Range("H8").Value = "Madrigal"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 140
$ws.Range("K8").Value = 15000
$ws.Range("L8").Value = 16000
$ws.Range("M8").Value = 15429
$ws.Range("N8").Value = "$/caja 40 unidades"
$ws.Range("O8").Value = "Provincia de Limarí"
$ws.Range("P8").Value = 386
$ws.Range("Q8").Value = 40
$ws.Range("R8").Value = "Hortaliza"
